$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: ReLU, camadas=2, neuronios=10, tempo, MSE
$ws.Range("A22").Value = "ReLU"
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 10
$ws.Range("D22").Value = [double]"1.208615E-2"
$ws.Range("E22").Value = [double]"2.8101820000000002"

# Row 23: ReLU, camadas=2, neuronios=15, tempo, MSE
$ws.Range("A23").Value = "ReLU"
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 15
$ws.Range("D23").Value = [double]"3.5551069999999997E-2"
$ws.Range("E23").Value = [double]"2.810168"

# Row 24: ReLU, camadas=3, neuronios=5 (no time/MSE result -> highlighted error row)
$ws.Range("A24").Value = "ReLU"
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 5
$ws.Range("A24:F24").Interior.Color = 255

# Row 25: ReLU, camadas=3, neuronios=10, tempo, MSE
$ws.Range("A25").Value = "ReLU"
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = [double]"3.2114980000000001E-2"
$ws.Range("E25").Value = [double]"2.810168"

# Row 26: ReLU, camadas=3, neuronios=15, tempo, MSE
$ws.Range("A26").Value = "ReLU"
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = [double]"5.3390029999999998E-2"
$ws.Range("E26").Value = [double]"2.810168"

$ws.Range("J17").Select() | Out-Null
